# Applies the diff: insert two new weekly records for Betarraga (rows 389-390)
# which pushes the existing data block down by two rows (old 389-482 -> new
# 391-484), extending the used range from A1:R482 to A1:R484.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 389. Excel will shift rows
# 389:482 down to 391:484, carrying their existing values/formatting
# (including the D-column date style) with them.
$ws.Rows("389:390").Insert()

# Row 391 is now an exact copy of the old row 389 ("Primera" template row),
# and row 392 is a copy of the old row 390 ("Segunda" template row). Copy
# those two rows into the freshly inserted 389/390 so every non-changing
# column (A,B,C,E,F,G,H,I,N,O,Q,R) is populated correctly, then overwrite
# the columns that actually change per the diff (D,J,K,L,M,P).
$ws.Range("A391:R391").Copy()
$ws.Range("A389:R389").PasteSpecial()
$ws.Range("A392:R392").Copy()
$ws.Range("A390:R390").PasteSpecial()
$excel.CutCopyMode = 0

# New "Primera" row (389)
$ws.Range("D389").Value = 44543
$ws.Range("J389").Value = 3100
$ws.Range("K389").Value = 400
$ws.Range("L389").Value = 450
$ws.Range("M389").Value = 426
$ws.Range("P389").Value = 106

# New "Segunda" row (390)
$ws.Range("D390").Value = 44543
$ws.Range("J390").Value = 1800
$ws.Range("K390").Value = 300
$ws.Range("L390").Value = 300
$ws.Range("M390").Value = 300
$ws.Range("P390").Value = 75
